$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Auditorias a procesos" section (rows 6-12): fill in the empty
#     "Estrategia de seguimiento" column (C) ---
$ws.Range("C6").Value = "Bimestral"
$ws.Range("C7").Value = "Al finalizar etapa de ventas"
$ws.Range("C8").Value = "Al finalizar etapa de planeación"
$ws.Range("C9").Value = "Al finalizar etapa de cierre"
$ws.Range("C10").Value = "Al finalizar etapa de cierre"
$ws.Range("C11").Value = "Mensual"
$ws.Range("C12").Value = "Semestral"

# --- "Auditorias a productos" section (rows 19-22) ---
$ws.Range("B19").Value = "Requerimientos"
$ws.Range("C19").Value = "Al finalizar etapa de ventas"

$ws.Range("B20").Value = "Estimación"
$ws.Range("C20").Value = "Al finalizar etapa de ventas"

$ws.Range("B21").Value = "Plan de proyecto"
$ws.Range("C21").Value = "Al finalizar etapa de planeación"

$ws.Range("C22").Value = "Semestral"

# --- Escalamiento section (rows 28-29): fill in the "Tiempo" column (C) ---
$ws.Range("C28").Value = "3 días habiles"
$ws.Range("C29").Value = "5 días habiles"

# Match the final selection left by the author
$ws.Range("C29").Select()
